$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (NCTId), shifting NCTId..intervention_type
# one column to the right, and making room for the new "statut_name" column.
$ws.Columns("C").Insert()

# Give the new header cell (C1) the same formatting as its neighbouring header
# cells (bold font + border), then set its value.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "statut_name"

# Fill the new column's data rows (2..last used row) with the constant value.
$lastRow = $ws.Cells($ws.Rows.Count, "A").End(-4162).Row
$ws.Range("C2:C" + $lastRow).Value = "pas de résultat ni de publication"
